# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2210"
#   "<header>_new" -> "<header>_FV2304"
# Then wrap the data range in a table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffixCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newSuffixCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $text = $cell.Value()
    $cell.Value = ($text -replace "_old$", "_FV2210")
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $text = $cell.Value()
    $cell.Value = ($text -replace "_new$", "_FV2304")
}

# Freeze the header row (row 1) so row 2 is the first scrollable row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# Wrap the used range in a table, matching the workbook's existing data extent.
$range = $ws.Range("A1:U74")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
